# The table on slide 6 ("SOURCES OF FINANCE") had its table style swapped
# to a different built-in PowerPoint table style.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{9414E99A-EC64-4228-8FBD-C17E7B2A91E0}")
    }
}
